# Modify train.py to read data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update SAVE_INTV (M2) from 100 to 200; N2 formula recalculates automatically.
$ws.Range("M2").Value = 200

# Update the active selection to M3, matching the saved view state.
$ws.Range("M3").Select()
